$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Ingreso")
$ws2 = $wb.Worksheets.Item("Gastos")

# --- Step 1: establish new shared-string insertion order (Iverson, Yandy, Mamao) ---
# This matches the order the strings first appear in the authored workbooks shared string table.
$ws1.Range("B573").Value = "Iverson"
$ws1.Range("B552").Value = "Yandy"   # was "Yandi" -> corrected to "Yandy"
$ws1.Range("B569").Value = "Mamao"

# --- Step 2: fill in the remaining new "Ingreso" (income) rows 554-576 ---
$ws1.Range("A554").Value = 45235
$ws1.Range("B554").Value = "Chamo"
$ws1.Range("C554").Value = 100
$ws1.Range("D554").Value = "Aporte"
$ws1.Range("A555").Value = 45235
$ws1.Range("B555").Value = "Wilkin"
$ws1.Range("C555").Value = 100
$ws1.Range("D555").Value = "Aporte"
$ws1.Range("A556").Value = 45235
$ws1.Range("B556").Value = "Invitados"
$ws1.Range("C556").Value = 100
$ws1.Range("D556").Value = "Aporte"
$ws1.Range("A557").Value = 45235
$ws1.Range("B557").Value = "Carlos"
$ws1.Range("C557").Value = 100
$ws1.Range("D557").Value = "Aporte"
$ws1.Range("A558").Value = 45235
$ws1.Range("B558").Value = "Rayder"
$ws1.Range("C558").Value = 100
$ws1.Range("D558").Value = "Aporte"
$ws1.Range("A559").Value = 45235
$ws1.Range("B559").Value = "Jeicol"
$ws1.Range("C559").Value = 100
$ws1.Range("D559").Value = "Aporte"
$ws1.Range("A560").Value = 45242
$ws1.Range("B560").Value = "Invitados"
$ws1.Range("C560").Value = 100
$ws1.Range("D560").Value = "Aporte"
$ws1.Range("A561").Value = 45242
$ws1.Range("B561").Value = "Julio"
$ws1.Range("C561").Value = 200
$ws1.Range("D561").Value = "Aporte"
$ws1.Range("A562").Value = 45242
$ws1.Range("B562").Value = "Melvin"
$ws1.Range("C562").Value = 100
$ws1.Range("D562").Value = "Aporte"
$ws1.Range("A563").Value = 45242
$ws1.Range("B563").Value = "Yandy"
$ws1.Range("C563").Value = 100
$ws1.Range("D563").Value = "Aporte"
$ws1.Range("A564").Value = 45242
$ws1.Range("B564").Value = "Rayder"
$ws1.Range("C564").Value = 100
$ws1.Range("D564").Value = "Aporte"
$ws1.Range("A565").Value = 45256
$ws1.Range("B565").Value = "Carlos"
$ws1.Range("C565").Value = 150
$ws1.Range("D565").Value = "Aporte"
$ws1.Range("A566").Value = 45256
$ws1.Range("B566").Value = "Yandy"
$ws1.Range("C566").Value = 100
$ws1.Range("D566").Value = "Aporte"
$ws1.Range("A567").Value = 45256
$ws1.Range("B567").Value = "Anuel"
$ws1.Range("C567").Value = 200
$ws1.Range("D567").Value = "Aporte"
$ws1.Range("A568").Value = 45256
$ws1.Range("B568").Value = "Kawai"
$ws1.Range("C568").Value = 100
$ws1.Range("D568").Value = "Aporte"
$ws1.Range("A569").Value = 45256
$ws1.Range("C569").Value = 100
$ws1.Range("D569").Value = "Aporte"
$ws1.Range("A570").Value = 45264
$ws1.Range("B570").Value = "Julio"
$ws1.Range("C570").Value = 200
$ws1.Range("D570").Value = "Aporte"
$ws1.Range("A571").Value = 45264
$ws1.Range("B571").Value = "Yeyo"
$ws1.Range("C571").Value = 100
$ws1.Range("D571").Value = "Aporte"
$ws1.Range("A572").Value = 45264
$ws1.Range("B572").Value = "Wilkin"
$ws1.Range("C572").Value = 100
$ws1.Range("D572").Value = "Aporte"
$ws1.Range("A573").Value = 45264
$ws1.Range("C573").Value = 10
$ws1.Range("D573").Value = "Aporte"
$ws1.Range("A574").Value = 45256
$ws1.Range("B574").Value = "Orlando"
$ws1.Range("C574").Value = 500
$ws1.Range("D574").Value = "Aporte"
$ws1.Range("A575").Value = 45259
$ws1.Range("B575").Value = "Joel"
$ws1.Range("C575").Value = 300
$ws1.Range("D575").Value = "Aporte"
$ws1.Range("A576").Value = 45264
$ws1.Range("B576").Value = "Johan"
$ws1.Range("C576").Value = 750
$ws1.Range("D576").Value = "Aporte"

# --- Step 3: fill in the new "Gastos" (expenses) rows 67-74 ---
$ws2.Range("A67").Value = 45221
$ws2.Range("B67").Value = "Arbitro, agua y hielo"
$ws2.Range("A68").Value = 45228
$ws2.Range("B68").Value = "Arbitro, agua y hielo"
$ws2.Range("C67").Formula = "=800+140"

# Rows 68-74 share one Arbitro formula (800+140) via a fill across the range,
# matching Excels native shared-formula serialisation.
$ws2.Range("A69").Value = 45235
$ws2.Range("B69").Value = "Arbitro, agua y hielo"
$ws2.Range("A70").Value = 45242
$ws2.Range("B70").Value = "Arbitro, agua y hielo"
$ws2.Range("A71").Value = 45249
$ws2.Range("B71").Value = "Arbitro, agua y hielo"
$ws2.Range("A72").Value = 45256
$ws2.Range("B72").Value = "Arbitro, agua y hielo"
$ws2.Range("A73").Value = 45263
$ws2.Range("B73").Value = "Arbitro, agua y hielo"
$ws2.Range("A74").Value = 45270
$ws2.Range("B74").Value = "Arbitro, agua y hielo"
$ws2.Range("C68:C74").Formula = "=800+140"

# --- Step 4: leave the selection where the author left it after the update ---
$ws1.Activate() | Out-Null
$ws1.Range("B563").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("C67").Select() | Out-Null
$ws1.Activate() | Out-Null
